$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.494.32'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '1.877.96'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7153'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.78'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07967'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3108'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.29'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.23%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08283'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7298'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.869.11'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.284'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.16'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('D16').Value = '29.485.01'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.930'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '245.68'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.24%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007883'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').Value = '2.119.06'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.956'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +6.17%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1611'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +13.80%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '163.94'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.062'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.32'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('E29').Value = '  -2.93%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.499'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.390'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.114'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05269'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('E34').Value = '  +2.31%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.199'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7276'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.97%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01869'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '1.224.58'
$ws.Range('E39').Value = '  +6.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.716'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9118'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '73.89'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.130'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.95%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.21'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').Value = '2.014.07'
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('E48').Value = '  +3.84%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.934'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +10.14%  '
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.342'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.12%  '
